$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.986.90'
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").Value = '1.872.54'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5077'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3666'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07214'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8952'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.78'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07525'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.73%  '
$ws.Range("B13").Value = 'Litecoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '95.34'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.48%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.866.48'
$ws.Range("E14").Value = '  +0.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.246'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.24%  '
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.24'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9998'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").Value = '27.023.04'
$ws.Range("E20").Value = '  -0.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.028'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").Value = '2.095.86'
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.40'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.399'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.40'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.40%  '
$ws.Range("E26").Value = '  -3.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.090'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.36'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.721'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.740'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09165'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05112'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7510'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.973'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.77%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.160'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("E37").Value = '  +6.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.552'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5650'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.86%  '
$ws.Range("E40").Value = '  -1.88%  '
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.649'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '115.61'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.576'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1476'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4777'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.98%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.16'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.55%  '
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9997'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.572'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.98'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.24%  '
